$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-21 Wednesday" "2026-01-22 Thursday"

Replace-Text "176÷3=58, 2" "837÷3=279, 0"
Replace-Text "192÷9=21, 3" "801÷9=89, 0"
Replace-Text "952÷8=119, 0" "488÷6=81, 2"
Replace-Text "410÷8=51, 2" "145÷5=29, 0"
Replace-Text "604÷6=100, 4" "758÷6=126, 2"

Replace-Text "896÷8=112, 0" "134÷3=44, 2"
Replace-Text "961÷5=192, 1" "283÷9=31, 4"
Replace-Text "260÷3=86, 2" "772÷3=257, 1"
Replace-Text "934÷6=155, 4" "142÷5=28, 2"
Replace-Text "299÷5=59, 4" "794÷5=158, 4"

Replace-Text "211÷5=42, 1" "149÷7=21, 2"
Replace-Text "872÷2=436, 0" "270÷9=30, 0"
Replace-Text "134÷2=67, 0" "135÷8=16, 7"
Replace-Text "705÷5=141, 0" "373÷6=62, 1"
Replace-Text "100÷2=50, 0" "636÷9=70, 6"

Replace-Text "568÷4=142, 0" "677÷8=84, 5"
Replace-Text "961÷7=137, 2" "761÷4=190, 1"
Replace-Text "797÷5=159, 2" "274÷9=30, 4"
Replace-Text "222÷3=74, 0" "294÷4=73, 2"
Replace-Text "126÷7=18, 0" "313÷4=78, 1"

Replace-Text "395÷4=98, 3" "837÷2=418, 1"
Replace-Text "525÷3=175, 0" "709÷9=78, 7"
Replace-Text "334÷4=83, 2" "508÷2=254, 0"
Replace-Text "376÷3=125, 1" "841÷8=105, 1"
Replace-Text "498÷4=124, 2" "228÷5=45, 3"

Write-Output "Done"
